$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 163
$ws.Range("I2").Value = 461
$ws.Range("J2").Value = 1809
$ws.Range("K2").Value = 9
$ws.Range("L2").Value = 529
$ws.Range("M2").Value = 30
$ws.Range("N2").Value = 306
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 9
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = 26
$ws.Range("S2").Value = 207
$ws.Range("T2").Value = 322
$ws.Range("U2").Value = 26
$ws.Range("V2").Value = 2736
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 2792
$ws.Range("Y2").Value = 3
$ws.Range("Z2").Value = 58
$ws.Range("AA2").Value = 19
